$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.21011
$ws.Range("H2").Value = 0.6303299999999999
$ws.Range("M2").Value = 12.80871533333333
$ws.Range("N2").Value = 38.426146
$ws.Range("O2").Value = 0.1716721220213608
$ws.Range("P2").Value = 0.1716721220213608
$ws.Range("Q2").Value = 2.691239178686667
$ws.Range("R2").Value = 24.22115260818
$ws.Range("S2").Value = 0.1716721220213608
$ws.Range("T2").Value = 0.1716721220213608
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.21011
$ws.Range("H3").Value = 0.6303299999999999
$ws.Range("N3").Value = 65.552207
$ws.Range("O3").Value = 0.2928601395225403
$ws.Range("P3").Value = 0.2928601395225403
$ws.Range("Q3").Value = 4.591058070923332
$ws.Range("R3").Value = 41.31952263830999
$ws.Range("S3").Value = 0.2928601395225403
$ws.Range("T3").Value = 0.2928601395225403
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.21011
$ws.Range("H4").Value = 0.6303299999999999
$ws.Range("M4").Value = 10.846871
$ws.Range("N4").Value = 32.540613
$ws.Range("O4").Value = 0.1453779956383313
$ws.Range("P4").Value = 0.1453779956383313
$ws.Range("Q4").Value = 2.27903606581
$ws.Range("R4").Value = 20.51132459229
$ws.Range("S4").Value = 0.1453779956383313
$ws.Range("T4").Value = 0.1453779956383313
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.21011
$ws.Range("H5").Value = 0.6303299999999999
$ws.Range("M5").Value = 7.389532
$ws.Range("N5").Value = 22.168596
$ws.Range("O5").Value = 0.09904011496636306
$ws.Range("P5").Value = 0.09904011496636304
$ws.Range("Q5").Value = 1.55261456852
$ws.Range("R5").Value = 13.97353111668
$ws.Range("S5").Value = 0.09904011496636306
$ws.Range("T5").Value = 0.09904011496636304
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.21011
$ws.Range("H6").Value = 0.6303299999999999
$ws.Range("M6").Value = 7.190038666666666
$ws.Range("N6").Value = 21.570116
$ws.Range("O6").Value = 0.09636635393950015
$ws.Range("P6").Value = 0.09636635393950013
$ws.Range("Q6").Value = 1.510699024253333
$ws.Range("R6").Value = 13.59629121828
$ws.Range("S6").Value = 0.09636635393950015
$ws.Range("T6").Value = 0.09636635393950013
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.21011
$ws.Range("H7").Value = 0.6303299999999999
$ws.Range("M7").Value = 14.525612
$ws.Range("N7").Value = 43.576836
$ws.Range("O7").Value = 0.1946832739119044
$ws.Range("P7").Value = 0.1946832739119044
$ws.Range("Q7").Value = 3.05197633732
$ws.Range("R7").Value = 27.46778703588
$ws.Range("S7").Value = 0.1946832739119044
$ws.Range("T7").Value = 0.1946832739119044
